$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column F (fastqFileName column)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 6).End(-4162).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $val = $cell.Value2
    if ($val -ne $null -and $val -is [string] -and $val.Contains("fastq_reformat/")) {
        $cell.Value = $val.Replace("fastq_reformat/", "")
    }
}
